# Update the "想去人数" (want-to-go count) figures that changed between this
# generated-output commit and the previous one. Two sheets are affected:
#   - "展览"   (sheet1) holds the primary rows for each event
#   - "全部类型" (sheet4) mirrors the same events (rows offset by +2)

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsAll     = $wb.Worksheets.Item("全部类型")

# Sheet "展览": cell -> new value
$exhibitUpdates = @{
    "F8"  = 1629
    "F9"  = 7433
    "F11" = 7625
    "F15" = 6132
    "F16" = 3251
    "F25" = 283
    "F35" = 1453
    "F40" = 154
    "F41" = 241
    "F44" = 477
    "F47" = 521
}

foreach ($cell in $exhibitUpdates.Keys) {
    $wsExhibit.Range($cell).Value = $exhibitUpdates[$cell]
}

# Sheet "全部类型": cell -> new value
$allUpdates = @{
    "F10" = 1629
    "F13" = 7433
    "F14" = 7625
    "F17" = 6132
    "F18" = 3251
    "F28" = 283
    "F39" = 1453
    "F44" = 241
    "F46" = 477
    "F48" = 521
}

foreach ($cell in $allUpdates.Keys) {
    $wsAll.Range($cell).Value = $allUpdates[$cell]
}
